$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the "Transposed" table (A9:B24) ascending by the
# "count of distinct values" column (B) instead of descending.
$ws.Range("A9:B24").Sort($ws.Range("B9"), 1)

# Move the selection (no scrolled topLeftCell this time)
$ws.Range("D20").Select() | Out-Null
